# Update EUR->ARS rate: 2025-10-06T15:20:01Z
# Appends a new data row (row 61) to the sheet with the latest quotation
# (date, time, and the EUR->ARS quote string), mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 61

# Column A holds a date-shaped string ("2025-10-06"). A leading apostrophe
# forces Excel to keep it as literal text instead of auto-converting it to
# a date serial number, matching how the other rows store plain text.
$ws.Range("A$newRow").Value = "'2025-10-06"
$ws.Range("B$newRow").Value = "15:20:01"
$ws.Range("C$newRow").Value = "1.00 EUR = 1,782.5648"
